{"js": "// Resume: rebuild resume items\n// Update the intro paragraph: \"3 years\" -> \"5 years\" and expand the\n// tooling sentence from \"the Laravel framework\" to a fuller tech list.\n\nconst body = context.document.body;\n\n// 1) \"... with 3 years commercial development experience\" -> \"... 5 years ...\"\nconst yearsResults = body.search(\"3 years commercial development experience\", { matchCase: true });\nyearsResults.load(\"items\");\nawait context.sync();\n\nif (yearsResults.items.length > 0) {\n  yearsResults.items[0].insertText(\n    \"5 years commercial development experience\",\n    Word.InsertLocation.replace\n  );\n}\n\n// 2) \"working with the Laravel framework.\" -> \"working with PHP (Laravel) and tools like Ansible, Docker and Terraform.\"\nconst frameworkResults = body.search(\"working with the Laravel framework.\", { matchCase: true });\nframeworkResults.load(\"items\");\nawait context.sync();\n\nif (frameworkResults.items.length > 0) {\n  frameworkResults.items[0].insertText(\n    \"working with PHP (Laravel) and tools like Ansible, Docker and Terraform.\",\n    Word.InsertLocation.replace\n  );\n}\n\nawait context.sync();\n", "ps1": "# Resume: rebuild resume items\n# Update the intro paragraph: \"3 years\" -> \"5 years\" and expand the\n# tooling sentence from \"the Laravel framework\" to a fuller tech list.\n\n$d = $word.ActiveDocument\n\n# 1) \"... with 3 years commercial development experience\" -> \"... 5 years ...\"\n$find1 = $d.Content.Find\n$find1.Text = \"3 years commercial development experience\"\n$find1.Replacement.Text = \"5 years commercial development experience\"\n$find1.Execute([ref]$find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find1.Replacement.Text, 2)\n\n# 2) \"working with the Laravel framework.\" -> \"working with PHP (Laravel) and tools like Ansible, Docker and Terraform.\"\n$find2 = $d.Content.Find\n$find2.Text = \"working with the Laravel framework.\"\n$find2.Replacement.Text = \"working with PHP (Laravel) and tools like Ansible, Docker and Terraform.\"\n$find2.Execute([ref]$find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find2.Replacement.Text, 2)\n"}
